# Refresh computed market-price / profit columns (H:N) on the per-job Leve
# profit sheets, as produced by the scheduled market-data runner.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: "Stuck in the Moment" (Leve Item ID 5505)
$ws.Range("H40").Value = 1806.5714
$ws.Range("I40").Value = 1776.8889
$ws.Range("J40").Value = 1860
$ws.Range("K40").Value = 1776.8889
$ws.Range("L40").Value = 1860
$ws.Range("M40").Value = -1601.8889
$ws.Range("N40").Value = -2210
# Row 64: "Forged from the Void" (Leve Item ID 5506)
$ws.Range("H64").Value = 3191.1
$ws.Range("I64").Value = 3237.5
$ws.Range("J64").Value = 3005.5
$ws.Range("K64").Value = 3237.5
$ws.Range("L64").Value = 3005.5
$ws.Range("M64").Value = -2989.5
$ws.Range("N64").Value = -3501.5
# Row 67: "Dodging the Draft (L)" (Leve Item ID 5506)
$ws.Range("H67").Value = 3191.1
$ws.Range("I67").Value = 3237.5
$ws.Range("J67").Value = 3005.5
$ws.Range("K67").Value = 3237.5
$ws.Range("L67").Value = 3005.5
$ws.Range("M67").Value = -2379.5
$ws.Range("N67").Value = -4721.5
# Row 70: "Consecrating Congregation" (Leve Item ID 12604)
$ws.Range("H70").Value = 1907.5264
$ws.Range("I70").Value = 859.1667
$ws.Range("J70").Value = 2391.3845
$ws.Range("K70").Value = 2577.5001
$ws.Range("L70").Value = 7174.1535
$ws.Range("M70").Value = -2307.5001
$ws.Range("N70").Value = -7714.1535
# Row 73: "Curbing the Contagion (L)" (Leve Item ID 12604)
$ws.Range("H73").Value = 1907.5264
$ws.Range("I73").Value = 859.1667
$ws.Range("J73").Value = 2391.3845
$ws.Range("K73").Value = 2577.5001
$ws.Range("L73").Value = 7174.1535
$ws.Range("M73").Value = -1641.5001
$ws.Range("N73").Value = -9046.1535
# Row 74: "Adhesive of Antipathy" (Leve Item ID 5507)
$ws.Range("H74").Value = 3496.6191
$ws.Range("I74").Value = 3045.5625
$ws.Range("K74").Value = 3045.5625
$ws.Range("M74").Value = -2109.5625
# Row 76: "Warding Off Temptation" (Leve Item ID 12602)
$ws.Range("H76").Value = 28574226
$ws.Range("I76").Value = 32260862
$ws.Range("J76").Value = 2800
$ws.Range("K76").Value = 32260862
$ws.Range("L76").Value = 2800
$ws.Range("M76").Value = -32260547
$ws.Range("N76").Value = -3430
# Row 77: "It's Gonna Grow Back (L)" (Leve Item ID 5507)
$ws.Range("H77").Value = 3496.6191
$ws.Range("I77").Value = 3045.5625
$ws.Range("K77").Value = 15227.8125
$ws.Range("M77").Value = -10547.8125
# Row 79: "The Garden of Arcane Delights (L)" (Leve Item ID 12602)
$ws.Range("H79").Value = 28574226
$ws.Range("I79").Value = 32260862
$ws.Range("J79").Value = 2800
$ws.Range("K79").Value = 32260862
$ws.Range("L79").Value = 2800
$ws.Range("M79").Value = -32259770
$ws.Range("N79").Value = -4984
# Row 139: "Something Salty and Ceremonial" (Leve Item ID 42306)
$ws.Range("H139").Value = 337966.66
$ws.Range("J139").Value = 337966.66
$ws.Range("L139").Value = 337966.66
$ws.Range("N139").Value = -348246.66

$ws = $wb.Worksheets.Item("ARM")
# Row 2: "Ain't Got No Ingots" (Leve Item ID 27713)
$ws.Range("H2").Value = 14838.6
$ws.Range("I2").Value = 19269.105
$ws.Range("K2").Value = 19269.105
$ws.Range("M2").Value = -19156.105
# Row 60: "Booty Call" (Leve Item ID 3883)
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
# Row 116: "No Scope" (Leve Item ID 27713)
$ws.Range("H116").Value = 14838.6
$ws.Range("I116").Value = 19269.105
$ws.Range("K116").Value = 19269.105
$ws.Range("M116").Value = -16975.105

$ws = $wb.Worksheets.Item("BSM")
# Row 3: "Hells Bells" (Leve Item ID 27713)
$ws.Range("H3").Value = 14838.6
$ws.Range("I3").Value = 19269.105
$ws.Range("K3").Value = 19269.105
$ws.Range("M3").Value = -19155.105
# Row 86: "Through Thick and Thin" (Leve Item ID 12526)
$ws.Range("H86").Value = 1865.847
$ws.Range("I86").Value = 1904.7051
$ws.Range("J86").Value = 1432.8572
$ws.Range("K86").Value = 1904.7051
$ws.Range("L86").Value = 1432.8572
$ws.Range("M86").Value = -781.7050999999999
$ws.Range("N86").Value = -3678.8572
# Row 89: "Piercing Eyes Deserve Piercing Shafts (L)" (Leve Item ID 12526)
$ws.Range("H89").Value = 1865.847
$ws.Range("I89").Value = 1904.7051
$ws.Range("J89").Value = 1432.8572
$ws.Range("K89").Value = 9523.5255
$ws.Range("L89").Value = 7164.286
$ws.Range("M89").Value = -3907.5255
$ws.Range("N89").Value = -18396.286
# Row 105: "Ingot to Wing It" (Leve Item ID 19947)
$ws.Range("H105").Value = 1928.5714
$ws.Range("I105").Value = 1928.5714
$ws.Range("K105").Value = 1928.5714
$ws.Range("M105").Value = -181.5714
# Row 134: "Ruthenium Supremium" (Leve Item ID 43998)
$ws.Range("H134").Value = 14671541
$ws.Range("I134").Value = 22477610
$ws.Range("J134").Value = 2702235.2
$ws.Range("K134").Value = 67432830
$ws.Range("L134").Value = 8106705.600000001
$ws.Range("M134").Value = -67430295
$ws.Range("N134").Value = -8111775.600000001

$ws = $wb.Worksheets.Item("CRP")
# Row 62: "Splinter in the Sewers" (Leve Item ID 12580)
$ws.Range("H62").Value = 2310
$ws.Range("I62").Value = 2260.625
$ws.Range("J62").Value = 2375.8333
$ws.Range("K62").Value = 2260.625
$ws.Range("L62").Value = 2375.8333
$ws.Range("M62").Value = -1636.625
$ws.Range("N62").Value = -3623.8333
# Row 65: "The Lumber of Their Discontent (L)" (Leve Item ID 12580)
$ws.Range("H65").Value = 2310
$ws.Range("I65").Value = 2260.625
$ws.Range("J65").Value = 2375.8333
$ws.Range("K65").Value = 11303.125
$ws.Range("L65").Value = 11879.1665
$ws.Range("M65").Value = -8183.125
$ws.Range("N65").Value = -18119.1665
# Row 134: "Wood You Be Quiet" (Leve Item ID 44020)
$ws.Range("H134").Value = 1606800.8
$ws.Range("I134").Value = 7221.6665
$ws.Range("J134").Value = 5720004
$ws.Range("K134").Value = 21664.9995
$ws.Range("L134").Value = 17160012
$ws.Range("M134").Value = -19129.9995
$ws.Range("N134").Value = -17165082
# Row 141: "No Greater Treasure" (Leve Item ID 43345)
$ws.Range("H141").Value = 136797.88
$ws.Range("J141").Value = 181598.33
$ws.Range("L141").Value = 181598.33
$ws.Range("N141").Value = -191958.33

$ws = $wb.Worksheets.Item("CUL")
# Row 113: "Can't Eat Just One" (Leve Item ID 27843)
$ws.Range("H113").Value = 1259.8679
$ws.Range("I113").Value = 1025.6957
$ws.Range("J113").Value = 1439.4
$ws.Range("K113").Value = 3077.0871
$ws.Range("L113").Value = 4318.200000000001
$ws.Range("M113").Value = -907.0870999999997
$ws.Range("N113").Value = -8658.200000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 70: "Sky Is the Limit" (Leve Item ID 14146)
$ws.Range("H70").Value = 4030604.5
$ws.Range("I70").Value = 1606561.4
$ws.Range("K70").Value = 1606561.4
$ws.Range("M70").Value = -1606291.4
# Row 73: "Hulls of Broken Dreams (L)" (Leve Item ID 14146)
$ws.Range("H73").Value = 4030604.5
$ws.Range("I73").Value = 1606561.4
$ws.Range("K73").Value = 1606561.4
$ws.Range("M73").Value = -1605625.4
# Row 80: "Needs More Prayerbell" (Leve Item ID 12521)
$ws.Range("H80").Value = 10920.904
$ws.Range("I80").Value = 4742.0713
$ws.Range("J80").Value = 23278.572
$ws.Range("K80").Value = 4742.0713
$ws.Range("L80").Value = 23278.572
$ws.Range("M80").Value = -3744.0713
$ws.Range("N80").Value = -25274.572
# Row 83: "With a Noise That Reaches Heaven (L)" (Leve Item ID 12521)
$ws.Range("H83").Value = 10920.904
$ws.Range("I83").Value = 4742.0713
$ws.Range("J83").Value = 23278.572
$ws.Range("K83").Value = 23710.3565
$ws.Range("L83").Value = 116392.86
$ws.Range("M83").Value = -18718.3565
$ws.Range("N83").Value = -126376.86
# Row 113: "Copious Crystal Cannons" (Leve Item ID 27710)
$ws.Range("H113").Value = 18334.2
$ws.Range("I113").Value = 1925.8889
$ws.Range("J113").Value = 42946.668
$ws.Range("K113").Value = 1925.8889
$ws.Range("L113").Value = 42946.668
$ws.Range("M113").Value = 244.1111000000001
$ws.Range("N113").Value = -47286.668

$ws = $wb.Worksheets.Item("LTW")
# Row 16: "Saddle Sore" (Leve Item ID 5289)
$ws.Range("H16").Value = 42383.793
$ws.Range("I16").Value = 53126.844
$ws.Range("J16").Value = 1560.2
$ws.Range("K16").Value = 53126.844
$ws.Range("L16").Value = 1560.2
$ws.Range("M16").Value = -52956.844
$ws.Range("N16").Value = -1900.2
